# Auto-generated edit script applying the diff changes to resum_diari_meteocat.xlsx
# (Update automatic: dades i banners [2026-02-24 04:49])
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-24 04:48:12'
$ws.Range('E3').Value = '2026-02-24 04:48:14'
$ws.Range('O3').Value = '2.8 °C'
$ws.Range('E4').Value = '2026-02-24 04:48:16'
$ws.Range('J4').Value = '1022.6 hPa'
$ws.Range('O4').Value = '6.9 °C'
$ws.Range('E5').Value = '2026-02-24 04:48:19'
$ws.Range('H5').Formula = "'36%"
$ws.Range('O5').Value = '3.8 °C'
$ws.Range('E6').Value = '2026-02-24 04:48:21'
$ws.Range('H6').Formula = "'81%"
$ws.Range('J6').Value = '1022.1 hPa'
$ws.Range('N6').Value = '8.2 °C 4:09 TU'
$ws.Range('O6').Value = '9.7 °C'
$ws.Range('E7').Value = '2026-02-24 04:48:23'
$ws.Range('L7').Value = '5.4 km/h - 289º 4:13 TU'
$ws.Range('N7').Value = '11.5 °C 4:00 TU'
$ws.Range('O7').Value = '12.3 °C'
$ws.Range('E8').Value = '2026-02-24 04:48:25'
$ws.Range('H8').Formula = "'46%"
$ws.Range('J8').Value = '1021.9 hPa'
$ws.Range('E9').Value = '2026-02-24 04:48:27'
$ws.Range('O9').Value = '5.7 °C'
$ws.Range('E10').Value = '2026-02-24 04:48:30'
$ws.Range('L10').Value = '5.4 km/h - 77º 4:21 TU'
$ws.Range('N10').Value = '3.9 °C 4:15 TU'
$ws.Range('O10').Value = '4.6 °C'
$ws.Range('E11').Value = '2026-02-24 04:48:32'
$ws.Range('H11').Formula = "'93%"
$ws.Range('N11').Value = '1.8 °C 4:28 TU'
$ws.Range('O11').Value = '2.8 °C'
$ws.Range('E12').Value = '2026-02-24 04:48:34'
$ws.Range('N12').Value = '4.0 °C 4:24 TU'
$ws.Range('O12').Value = '6.3 °C'
$ws.Range('E13').Value = '2026-02-24 04:48:37'
$ws.Range('J13').Value = '1029.9 hPa'
$ws.Range('N13').Value = '-3.2 °C 4:18 TU'
$ws.Range('O13').Value = '-1.1 °C'
$ws.Range('E14').Value = '2026-02-24 04:48:39'
$ws.Range('N14').Value = '8.2 °C 4:18 TU'
$ws.Range('O14').Value = '9.1 °C'
$ws.Range('E15').Value = '2026-02-24 04:48:41'
$ws.Range('H15').Formula = "'91%"
$ws.Range('E16').Value = '2026-02-24 04:48:43'
$ws.Range('H16').Formula = "'21%"
$ws.Range('K16').Value = '-0.1 MJ/m2'
$ws.Range('N16').Value = '3.3 °C 4:12 TU'
$ws.Range('E17').Value = '2026-02-24 04:48:45'
$ws.Range('E18').Value = '2026-02-24 04:48:48'
$ws.Range('J18').Value = '1022.8 hPa'
$ws.Range('N18').Value = '2.1 °C 4:15 TU'
$ws.Range('O18').Value = '3.3 °C'
$ws.Range('E19').Value = '2026-02-24 04:48:50'
$ws.Range('H19').Formula = "'70%"
$ws.Range('O19').Value = '9.1 °C'
$ws.Range('E20').Value = '2026-02-24 04:48:52'
$ws.Range('H20').Formula = "'44%"
$ws.Range('L20').Value = '18.4 km/h - 259º 4:17 TU'
$ws.Range('N20').Value = '0.3 °C 4:29 TU'
$ws.Range('O20').Value = '1.1 °C'
$ws.Range('E21').Value = '2026-02-24 04:48:55'
$ws.Range('N21').Value = '2.3 °C 4:29 TU'
$ws.Range('O21').Value = '3.8 °C'
$ws.Range('E22').Value = '2026-02-24 04:48:57'
$ws.Range('H22').Formula = "'19%"
$ws.Range('E23').Value = '2026-02-24 04:48:59'
$ws.Range('H23').Formula = "'29%"
$ws.Range('L23').Value = '10.4 km/h - 51º 4:14 TU'
$ws.Range('E24').Value = '2026-02-24 04:49:02'
$ws.Range('N24').Value = '1.6 °C 4:26 TU'
$ws.Range('O24').Value = '3.7 °C'
$ws.Range('E25').Value = '2026-02-24 04:49:04'
$ws.Range('H25').Formula = "'36%"
$ws.Range('O25').Value = '4.7 °C'
$ws.Range('E26').Value = '2026-02-24 04:49:06'
$ws.Range('H26').Formula = "'55%"
$ws.Range('J26').Value = '1022.8 hPa'
$ws.Range('O26').Value = '7.6 °C'
$ws.Range('E27').Value = '2026-02-24 04:49:09'
$ws.Range('H27').Formula = "'40%"
$ws.Range('O27').Value = '4.3 °C'
$ws.Range('E28').Value = '2026-02-24 04:49:11'
$ws.Range('J28').Value = '1024.1 hPa'
$ws.Range('O28').Value = '3.9 °C'
$ws.Range('E29').Value = '2026-02-24 04:49:13'
$ws.Range('N29').Value = '3.9 °C 4:08 TU'
$ws.Range('O29').Value = '5.4 °C'
$ws.Range('E30').Value = '2026-02-24 04:49:15'
$ws.Range('J30').Value = '1022.1 hPa'
$ws.Range('N30').Value = '8.0 °C 4:11 TU'
$ws.Range('O30').Value = '9.4 °C'
$ws.Range('E31').Value = '2026-02-24 04:49:18'
$ws.Range('J31').Value = '1021.2 hPa'
$ws.Range('O31').Value = '14.9 °C'
$ws.Range('E32').Value = '2026-02-24 04:49:20'
$ws.Range('N32').Value = '-4.1 °C 4:29 TU'
$ws.Range('O32').Value = '-3.0 °C'
$ws.Range('E33').Value = '2026-02-24 04:49:22'
$ws.Range('H33').Formula = "'68%"
$ws.Range('N33').Value = '0.8 °C 4:29 TU'
$ws.Range('O33').Value = '2.4 °C'
$ws.Range('E34').Value = '2026-02-24 04:49:24'
$ws.Range('H34').Formula = "'56%"
$ws.Range('M34').Value = '5.1 °C 4:17 TU'
$ws.Range('O34').Value = '2.6 °C'
$ws.Range('E35').Value = '2026-02-24 04:49:27'
$ws.Range('J35').Value = '1024.6 hPa'
$ws.Range('E36').Value = '2026-02-24 04:49:29'
$ws.Range('J36').Value = '1022.0 hPa'
$ws.Range('N36').Value = '6.9 °C 4:16 TU'
$ws.Range('O36').Value = '8.3 °C'
$ws.Range('E37').Value = '2026-02-24 04:49:31'
$ws.Range('J37').Value = '1027.7 hPa'
$ws.Range('N37').Value = '-0.3 °C 4:24 TU'
$ws.Range('O37').Value = '0.9 °C'
$ws.Range('E38').Value = '2026-02-24 04:49:34'
$ws.Range('H38').Formula = "'83%"
$ws.Range('E39').Value = '2026-02-24 04:49:36'
$ws.Range('H39').Formula = "'35%"
$ws.Range('K39').Value = '-0.1 MJ/m2'
$ws.Range('O39').Value = '4.8 °C'
$ws.Range('E40').Value = '2026-02-24 04:49:39'
$ws.Range('J40').Value = '1027.9 hPa'
$ws.Range('N40').Value = '0.3 °C 4:02 TU'
$ws.Range('O40').Value = '1.4 °C'
$ws.Range('E41').Value = '2026-02-24 04:49:41'
$ws.Range('H41').Formula = "'79%"
$ws.Range('J41').Value = '1022.4 hPa'
$ws.Range('K41').Value = '-0.1 MJ/m2'
$ws.Range('O41').Value = '7.2 °C'
$ws.Range('E42').Value = '2026-02-24 04:49:43'
$ws.Range('H42').Formula = "'96%"
$ws.Range('N42').Value = '5.4 °C 4:16 TU'
$ws.Range('O42').Value = '7.0 °C'
$ws.Range('E43').Value = '2026-02-24 04:49:46'
$ws.Range('O43').Value = '4.4 °C'
$ws.Range('E44').Value = '2026-02-24 04:49:48'
$ws.Range('L44').Value = '6.1 km/h - 26º 4:26 TU'
$ws.Range('O44').Value = '0.2 °C'
$ws.Range('E45').Value = '2026-02-24 04:49:50'
$ws.Range('H45').Formula = "'61%"
$ws.Range('N45').Value = '2.4 °C 4:09 TU'
$ws.Range('O45').Value = '4.5 °C'
$ws.Range('E46').Value = '2026-02-24 04:49:53'
$ws.Range('J46').Value = '1024.0 hPa'
$ws.Range('N46').Value = '1.0 °C 4:29 TU'
$ws.Range('O46').Value = '2.6 °C'
